# Generate Report for Handoff
# Updates status + handoff/generate timestamps for the zh-cn and de-de
# localization rows, and widens the status columns to fit the new,
# longer "Ready for handoff" text (mirrors Excel's own column autosize).

$wb = $excel.ActiveWorkbook

$newWidth = 16.3333333333333  # closest reachable ColumnWidth -> xml width ~17.22

# --- Overview sheet: one row summarizing both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"          # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"          # de-de status
$wsOverview.Range("G2").Value = "2016-08-20 22:45:52"        # Latest HO Xliff Generate Date
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"                # Status
$wsZh.Range("H2").Value = "2016-08-20 22:45:48"              # Latest Handoff Datetime
$wsZh.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"                # Status
$wsDe.Range("H2").Value = "2016-08-20 22:45:52"              # Latest Handoff Datetime
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
